# Update the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.109.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.363.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.43%  '

$ws.Range("E4").Value = '  +0.61%  '

$ws.Range("E5").Value = '  -1.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.05%  '

$ws.Range("E7").Value = '  +0.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.544'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.366.08'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0984'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.42%  '

$ws.Range("E11").Value = '  -0.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.81'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.29%  '

$ws.Range("E13").Value = '  +0.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.783.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.51%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.056.48'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.30%  '

$ws.Range("E17").Value = '  -2.57%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.306.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.49%  '

$ws.Range("E19").Value = '  -2.96%  '

$ws.Range("E20").Value = '  -2.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '307.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.77%  '

$ws.Range("E23").Value = '  +0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.368'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.34%  '

$ws.Range("E27").Value = '  -5.51%  '

$ws.Range("E28").Value = '  -4.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0711'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.33%  '

$ws.Range("E31").Value = '  -2.42%  '

$ws.Range("E32").Value = '  +0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.67%  '

$ws.Range("E34").Value = '  +0.60%  '

$ws.Range("E35").Value = '  -4.99%  '

$ws.Range("E36").Value = '  -2.45%  '

$ws.Range("E37").Value = '  -5.56%  '

$ws.Range("E38").Value = '  -4.46%  '

$ws.Range("E39").Value = '  -1.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.798'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.43%  '

$ws.Range("E41").Value = '  -5.91%  '

$ws.Range("E42").Value = '  -0.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '129.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.42%  '

$ws.Range("E44").Value = '  -6.75%  '

$ws.Range("E45").Value = '  -1.84%  '

$ws.Range("E46").Value = '  -1.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '238.23'
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0481'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0206'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.51%  '

$ws.Range("E50").Value = '  -1.64%  '

$ws.Range("E51").Value = '  -1.15%  '
